# Apply the cryptos-list refresh described by the commit diff.
# D-column cells whose new text looks like a plain number must be forced to
# Text format first, otherwise Excel auto-converts the literal string into a
# numeric value (losing formatting like trailing zeros, e.g. "15.50" -> 15.5).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numericLookingCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D11",
    "D12",
    "D13",
    "D14",
    "D15",
    "D16",
    "D18",
    "D19",
    "D20",
    "D22",
    "D23",
    "D24",
    "D25",
    "D27",
    "D28",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)

foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.377.41'
$ws.Range("E2").Value = '  -0.42%  '
$ws.Range("D3").Value = '1.720.79'
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '241.89'
$ws.Range("E5").Value = '  -1.98%  '
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").Value = '0.4882'
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("D8").Value = '0.2592'
$ws.Range("E8").Value = '  -2.75%  '
$ws.Range("D9").Value = '0.06181'
$ws.Range("E9").Value = '  -0.65%  '
$ws.Range("D10").Value = '1.725.99'
$ws.Range("E10").Value = '  -0.32%  '
$ws.Range("D11").Value = '0.06968'
$ws.Range("E11").Value = '  -1.20%  '
$ws.Range("D12").Value = '15.50'
$ws.Range("E12").Value = '  -1.25%  '
$ws.Range("D13").Value = '4.516'
$ws.Range("E13").Value = '  -2.27%  '
$ws.Range("D14").Value = '0.5969'
$ws.Range("E14").Value = '  -2.07%  '
$ws.Range("D15").Value = '77.04'
$ws.Range("D16").Value = '0.9994'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").Value = '26.372.92'
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").Value = '0.9993'
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("D19").Value = '0.000007156'
$ws.Range("E19").Value = '  -1.46%  '
$ws.Range("D20").Value = '11.31'
$ws.Range("E20").Value = '  -1.72%  '
$ws.Range("D21").Value = '1.947.93'
$ws.Range("E21").Value = '  -0.38%  '
$ws.Range("D22").Value = '4.438'
$ws.Range("E22").Value = '  -1.92%  '
$ws.Range("D23").Value = '8.477'
$ws.Range("E23").Value = '  -3.38%  '
$ws.Range("D24").Value = '5.068'
$ws.Range("E24").Value = '  -3.50%  '
$ws.Range("D25").Value = '137.80'
$ws.Range("E26").Value = '  -1.62%  '
$ws.Range("D27").Value = '1.399'
$ws.Range("E27").Value = '  -0.49%  '
$ws.Range("D28").Value = '106.25'
$ws.Range("E28").Value = '  -1.76%  '
$ws.Range("E29").Value = '  -3.13%  '
$ws.Range("D30").Value = '3.897'
$ws.Range("E30").Value = '  -1.98%  '
$ws.Range("D31").Value = '0.08008'
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("D32").Value = '3.652'
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("D33").Value = '0.04491'
$ws.Range("E33").Value = '  -1.63%  '
$ws.Range("D34").Value = '0.9986'
$ws.Range("D35").Value = '2.602'
$ws.Range("E35").Value = '  -0.46%  '
$ws.Range("D36").Value = '0.9949'
$ws.Range("E36").Value = '  -1.36%  '
$ws.Range("D37").Value = '0.6212'
$ws.Range("E37").Value = '  -2.65%  '
$ws.Range("D38").Value = '0.9183'
$ws.Range("E38").Value = '  +2.20%  '
$ws.Range("D39").Value = '1.961'
$ws.Range("E39").Value = '  -3.58%  '
$ws.Range("D40").Value = '2.380'
$ws.Range("E40").Value = '  -0.80%  '
$ws.Range("D41").Value = '0.9987'
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("D42").Value = '0.01476'
$ws.Range("E42").Value = '  -2.04%  '
$ws.Range("D43").Value = '100.24'
$ws.Range("E43").Value = '  -1.19%  '
$ws.Range("D44").Value = '5.430'
$ws.Range("E44").Value = '  -0.23%  '
$ws.Range("D45").Value = '0.3840'
$ws.Range("E45").Value = '  -1.24%  '
$ws.Range("D46").Value = '6.901'
$ws.Range("E46").Value = '  -0.69%  '
$ws.Range("D47").Value = '0.1162'
$ws.Range("E47").Value = '  -1.87%  '
$ws.Range("D48").Value = '0.05364'
$ws.Range("E48").Value = '  -0.44%  '
$ws.Range("D49").Value = '30.21'
$ws.Range("E49").Value = '  -1.13%  '
$ws.Range("D50").Value = '7.672'
$ws.Range("E50").Value = '  -1.65%  '
$ws.Range("D51").Value = '50.97'
$ws.Range("E51").Value = '  -0.75%  '

# Restore the default (General) style on the cells we forced to Text so the
# only observable change is the cell value, matching the source diff exactly.
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
